$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B35 / C35 to date values matching the style already used elsewhere (row 34)
$ws.Range("B34").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("C34").Copy()
$ws.Range("C35").PasteSpecial(-4122)
$ws.Range("B35").Value = 45516
$ws.Range("C35").Value = 45516
$excel.CutCopyMode = $false

# New formulas for P40, P42, P43, P44, P47, P49
$ws.Range("P40").Formula = "=60*60*100000*0.001/(0.001*0.001)"
$ws.Range("P42").Formula = "=0.45/(0.0254*0.0254)"
$ws.Range("P43").Formula = "=393149903724510"
$ws.Range("P44").Formula = "=60*60*6896/(100*0.00454609*0.00454609)"
$ws.Range("P47").Formula = "=1/0.000145"
$ws.Range("P49").Formula = "=0.453592*9.85/(0.0254*0.0254)"

# Update selection to match diff
$ws.Range("S48").Select()
